$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update ConceptScheme URI (B1)
$ws.Range("B1").Value = "http://ontology.deic.dk/cv/rock-n-roll/"

# Update PREFIX URI (C2)
$ws.Range("C2").Value = "http://vocab.deic.dk/cv/rock-n-roll/"

# Row 18: add a new test term
$ws.Range("A18").Value = "rock-n-roll:TestTerm"
$ws.Range("B18").Value = "TestTerm"
$ws.Range("E18").Value = "rock-n-roll:Deprecated"
